$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos price/volume refresh (GitHub Actions scrape), Tue Dec 5 08:58:43 UTC 2023.
#
# Columns D (Price) and E (Volume(1h)) store numeric-looking values as plain
# text (e.g. "41.592.89", "1.00", "  -0.32%  "). Excel's automatic type
# detection would otherwise coerce strings like "1.00" or "229.06" into real
# numbers when assigned via .Value, silently dropping formatting (e.g.
# trailing zeros) and changing the cell type. Each cell we touch below is
# first forced to Text format ("@") so the literal string is preserved,
# matching the existing (pre-edit) cell formatting used throughout the sheet.
#
# Rows 47/48 also swap identity: the scraper now lists FTXToken ahead of
# Maker, so what was row 47 (Maker) becomes FTXToken and vice versa, each
# carrying its own freshly scraped price/volume figures.
$changes = @(
    @{Addr='B47'; Value='FTXToken'},
    @{Addr='C47'; Value='https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'},
    @{Addr='D47'; Value='4.35'},
    @{Addr='E47'; Value='  -15.39%  '},
    @{Addr='B48'; Value='Maker'},
    @{Addr='C48'; Value='https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'},
    @{Addr='D48'; Value='1.459.66'},
    @{Addr='E48'; Value='  -3.49%  '},
    @{Addr='D2'; Value='41.592.89'},
    @{Addr='E2'; Value='  -0.32%  '},
    @{Addr='D3'; Value='2.198.92'},
    @{Addr='E3'; Value='  -2.61%  '},
    @{Addr='E4'; Value='  -0.02%  '},
    @{Addr='D5'; Value='229.06'},
    @{Addr='E5'; Value='  -1.78%  '},
    @{Addr='E6'; Value='  -4.20%  '},
    @{Addr='D7'; Value='60.01'},
    @{Addr='E7'; Value='  -6.73%  '},
    @{Addr='E9'; Value='  -3.39%  '},
    @{Addr='D10'; Value='56.80'},
    @{Addr='E10'; Value='  -5.99%  '},
    @{Addr='E11'; Value='  -2.88%  '},
    @{Addr='E12'; Value='  -2.09%  '},
    @{Addr='D13'; Value='2.522.23'},
    @{Addr='E13'; Value='  -2.73%  '},
    @{Addr='D14'; Value='15.33'},
    @{Addr='E14'; Value='  -5.64%  '},
    @{Addr='D15'; Value='22.13'},
    @{Addr='E15'; Value='  -2.50%  '},
    @{Addr='D16'; Value='5.64'},
    @{Addr='E16'; Value='  -0.90%  '},
    @{Addr='E17'; Value='  -4.92%  '},
    @{Addr='D18'; Value='2.212.61'},
    @{Addr='E18'; Value='  -2.08%  '},
    @{Addr='D19'; Value='41.472.59'},
    @{Addr='E19'; Value='  -0.24%  '},
    @{Addr='D20'; Value='71.81'},
    @{Addr='E20'; Value='  -3.02%  '},
    @{Addr='E21'; Value='  -4.13%  '},
    @{Addr='E22'; Value='  -3.13%  '},
    @{Addr='D23'; Value='241.78'},
    @{Addr='E23'; Value='  -4.74%  '},
    @{Addr='E24'; Value='  -0.11%  '},
    @{Addr='D25'; Value='2.34'},
    @{Addr='E25'; Value='  -2.79%  '},
    @{Addr='E26'; Value='  -2.16%  '},
    @{Addr='E27'; Value='  -3.18%  '},
    @{Addr='D28'; Value='168.59'},
    @{Addr='E28'; Value='  -1.78%  '},
    @{Addr='E29'; Value='  -7.14%  '},
    @{Addr='D30'; Value='1.45'},
    @{Addr='E30'; Value='  -0.41%  '},
    @{Addr='D31'; Value='19.67'},
    @{Addr='E31'; Value='  -4.27%  '},
    @{Addr='D32'; Value='2.57'},
    @{Addr='E32'; Value='  -8.92%  '},
    @{Addr='E33'; Value='  -4.33%  '},
    @{Addr='E34'; Value='  -2.75%  '},
    @{Addr='D35'; Value='4.61'},
    @{Addr='E35'; Value='  -3.62%  '},
    @{Addr='D36'; Value='0.0644'},
    @{Addr='E36'; Value='  +0.16%  '},
    @{Addr='D37'; Value='2.35'},
    @{Addr='E37'; Value='  -5.07%  '},
    @{Addr='E38'; Value='  -9.22%  '},
    @{Addr='D39'; Value='3.51'},
    @{Addr='E39'; Value='  -8.89%  '},
    @{Addr='E40'; Value='  -7.84%  '},
    @{Addr='D41'; Value='1.00'},
    @{Addr='E41'; Value='  +0.10%  '},
    @{Addr='E42'; Value='  -3.39%  '},
    @{Addr='D43'; Value='8.47'},
    @{Addr='E44'; Value='  -4.94%  '},
    @{Addr='E45'; Value='  -3.71%  '},
    @{Addr='D46'; Value='96.68'},
    @{Addr='D49'; Value='16.23'},
    @{Addr='E49'; Value='  -8.20%  '},
    @{Addr='E50'; Value='  -1.17%  '},
    @{Addr='E51'; Value='  -7.24%  '}
)

foreach ($item in $changes) {
    $cell = $ws.Range($item.Addr)
    $cell.NumberFormat = "@"
    $cell.Value = $item.Value
}
